$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
Write-Host "before:" $ws.AutoFilter.Range.Address()
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "Add Visitor Manually"
$ws.Range("C41").Value = "Shambhoo"
$ws.Range("D41").Value = "Pending"
$ws.Range("E41").Value = "Visitor"
Write-Host "after values:" $ws.AutoFilter.Range.Address()
